$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D3 to "PLEDGE" and copy the italic-style formatting used by D2/D4
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D3").Value = "PLEDGE"

# Set E4 to the new in-favor-of-whom value
$ws.Range("E4").Value = "в пользу кого 3"

# Update the selection to match the target state
$ws.Range("E3:E4").Select()
